$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 21666.334
$ws.Range("I13").Value = 10000
$ws.Range("J13").Value = 27499.5
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 27499.5
$ws.Range("M13").Value = -9831
$ws.Range("N13").Value = -27837.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5257.1
$ws.Range("I106").Value = 4696.5
$ws.Range("K106").Value = 4696.5
$ws.Range("M106").Value = -4065.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4990.375
$ws.Range("J112").Value = 5390.7144
$ws.Range("L112").Value = 16172.1432
$ws.Range("N112").Value = -18388.1432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3330.3333
$ws.Range("I113").Value = 3716.3333
$ws.Range("J113").Value = 2558.3333
$ws.Range("K113").Value = 3716.3333
$ws.Range("L113").Value = 2558.3333
$ws.Range("M113").Value = -462.3332999999998
$ws.Range("N113").Value = -9066.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2519.875
$ws.Range("I129").Value = 2532.5
$ws.Range("J129").Value = 2515.6667
$ws.Range("K129").Value = 7597.5
$ws.Range("L129").Value = 7547.000100000001
$ws.Range("M129").Value = -2597.5
$ws.Range("N129").Value = -17547.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2706.7693
$ws.Range("I132").Value = 2265.3901
$ws.Range("K132").Value = 6796.1703
$ws.Range("M132").Value = -4266.1703

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 6847.8184
$ws.Range("I135").Value = 2482.1052
$ws.Range("J135").Value = 34497.332
$ws.Range("K135").Value = 22338.9468
$ws.Range("L135").Value = 310475.988
$ws.Range("M135").Value = -19803.9468
$ws.Range("N135").Value = -315545.988

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3859.779
$ws.Range("I138").Value = 1871.5483
$ws.Range("J138").Value = 4822.828
$ws.Range("K138").Value = 5614.644899999999
$ws.Range("L138").Value = 14468.484
$ws.Range("M138").Value = -474.6448999999993
$ws.Range("N138").Value = -24748.484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10265.027
$ws.Range("J32").Value = 31910
$ws.Range("L32").Value = 31910
$ws.Range("N32").Value = -32484

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 14095.823
$ws.Range("I45").Value = 20122.637
$ws.Range("K45").Value = 20122.637
$ws.Range("M45").Value = -19745.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4863
$ws.Range("I61").Value = 4737.077
$ws.Range("J61").Value = 6500
$ws.Range("K61").Value = 4737.077
$ws.Range("L61").Value = 6500
$ws.Range("M61").Value = -4525.077
$ws.Range("N61").Value = -6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H81").Value = 144180.33
$ws.Range("I81").Value = 144180
$ws.Range("K81").Value = 144180
$ws.Range("M81").Value = -143182

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 27500
$ws.Range("I82").Value = 20000
$ws.Range("J82").Value = 35000
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 35000
$ws.Range("M82").Value = -19639
$ws.Range("N82").Value = -35722

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H84").Value = 144180.33
$ws.Range("I84").Value = 144180
$ws.Range("K84").Value = 432540
$ws.Range("M84").Value = -427548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H85").Value = 27500
$ws.Range("I85").Value = 20000
$ws.Range("J85").Value = 35000
$ws.Range("K85").Value = 20000
$ws.Range("L85").Value = 35000
$ws.Range("M85").Value = -18752
$ws.Range("N85").Value = -37496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4863
$ws.Range("I136").Value = 4737.077
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 14211.231
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -11661.231
$ws.Range("N136").Value = -24600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2397.5
$ws.Range("I105").Value = 2360.5454
$ws.Range("K105").Value = 2360.5454
$ws.Range("M105").Value = -613.5454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1997.325
$ws.Range("I134").Value = 1576.421
$ws.Range("K134").Value = 4729.263
$ws.Range("M134").Value = -2194.263

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2748.9375
$ws.Range("I31").Value = 2488.75
$ws.Range("J31").Value = 2905.05
$ws.Range("K31").Value = 2488.75
$ws.Range("L31").Value = 2905.05
$ws.Range("M31").Value = -2193.75
$ws.Range("N31").Value = -3495.05

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2748.9375
$ws.Range("I34").Value = 2488.75
$ws.Range("J34").Value = 2905.05
$ws.Range("K34").Value = 2488.75
$ws.Range("L34").Value = 2905.05
$ws.Range("M34").Value = -2286.75
$ws.Range("N34").Value = -3309.05

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1546.95
$ws.Range("I107").Value = 961.7143
$ws.Range("K107").Value = 961.7143
$ws.Range("M107").Value = 958.2857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 28000
$ws.Range("J112").Value = 28000
$ws.Range("L112").Value = 28000
$ws.Range("N112").Value = -30954

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 90945624
$ws.Range("J37").Value = 90945624
$ws.Range("L37").Value = 272836872
$ws.Range("N37").Value = -272837096

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7487.9443
$ws.Range("I56").Value = 7487.9443
$ws.Range("K56").Value = 7487.9443
$ws.Range("M56").Value = -6957.9443

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 1673.75
$ws.Range("I99").Value = 847.5
$ws.Range("K99").Value = 2542.5
$ws.Range("M99").Value = -296.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1078.4546
$ws.Range("I134").Value = 1078.4546
$ws.Range("K134").Value = 3235.3638
$ws.Range("M134").Value = 1834.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4086.2632
$ws.Range("J137").Value = 4504
$ws.Range("L137").Value = 13512
$ws.Range("N137").Value = -23712

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4818.625
$ws.Range("I22").Value = 274.5
$ws.Range("J22").Value = 6333.3335
$ws.Range("K22").Value = 274.5
$ws.Range("L22").Value = 6333.3335
$ws.Range("M22").Value = 254.5
$ws.Range("N22").Value = -7391.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3912.2173
$ws.Range("I122").Value = 2605.516
$ws.Range("K122").Value = 7816.548000000001
$ws.Range("M122").Value = -5366.548000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5992.385
$ws.Range("I126").Value = 4013.4285
$ws.Range("K126").Value = 12040.2855
$ws.Range("M126").Value = -9570.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13284
$ws.Range("I7").Value = 11942.23
$ws.Range("K7").Value = 11942.23
$ws.Range("M7").Value = -11830.23

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2577.6316
$ws.Range("I16").Value = 2211.7334
$ws.Range("J16").Value = 3949.75
$ws.Range("K16").Value = 2211.7334
$ws.Range("L16").Value = 3949.75
$ws.Range("M16").Value = -2041.7334
$ws.Range("N16").Value = -4289.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6627.0386
$ws.Range("I40").Value = 3049.1177
$ws.Range("J40").Value = 13385.333
$ws.Range("K40").Value = 3049.1177
$ws.Range("L40").Value = 13385.333
$ws.Range("M40").Value = -2913.1177
$ws.Range("N40").Value = -13657.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 216.71428
$ws.Range("I55").Value = 120
$ws.Range("K55").Value = 120
$ws.Range("M55").Value = 53

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3279
$ws.Range("I68").Value = 2966
$ws.Range("K68").Value = 2966
$ws.Range("M68").Value = -2217

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3279
$ws.Range("I71").Value = 2966
$ws.Range("K71").Value = 14830
$ws.Range("M71").Value = -11086

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 13284
$ws.Range("I126").Value = 11942.23
$ws.Range("K126").Value = 35826.69
$ws.Range("M126").Value = -33356.69

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 352389.66
$ws.Range("J4").Value = 675
$ws.Range("L4").Value = 675
$ws.Range("N4").Value = -901

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2261.3333
$ws.Range("I100").Value = 2261.3333
$ws.Range("K100").Value = 4522.6666
$ws.Range("M100").Value = -3981.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5126.0415
$ws.Range("I126").Value = 5301.5884
$ws.Range("J126").Value = 4699.7144
$ws.Range("K126").Value = 15904.7652
$ws.Range("L126").Value = 14099.1432
$ws.Range("M126").Value = -13434.7652
$ws.Range("N126").Value = -19039.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1558.4
$ws.Range("I136").Value = 1495
$ws.Range("J136").Value = 1639.091
$ws.Range("K136").Value = 4485
$ws.Range("L136").Value = 4917.272999999999
$ws.Range("M136").Value = -1935
$ws.Range("N136").Value = -10017.273
